$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$ws.Cells.Item(2,5).Value = "'752.725639062669"
$ws.Cells.Item(52,5).Value = "'757.592399660057"
$ws.Cells.Item(93,5).Value = "'1073.79218515661"
$ws.Cells.Item(94,5).Value = "'1170.12947504923"
$ws.Cells.Item(95,5).Value = "'1216.45975192823"
$ws.Cells.Item(96,5).Value = "'1219.25016212823"
$ws.Cells.Item(97,5).Value = "'1286.29796086775"
$ws.Cells.Item(98,5).Value = "'1385.46144533211"
$ws.Cells.Item(99,5).Value = "'1307.7706708874"
$ws.Cells.Item(100,5).Value = "'1382.31081619543"
$ws.Cells.Item(101,5).Value = "'1102.64694471882"
$ws.Cells.Item(102,5).Value = "'1469.82078286337"
$ws.Cells.Item(103,5).Value = "'1440.44329272521"
$ws.Cells.Item(104,5).Value = "'1433.48686075922"
$ws.Cells.Item(105,5).Value = "'1453.13020548845"
$ws.Cells.Item(106,5).Value = "'1402.92713690053"
$ws.Cells.Item(107,5).Value = "'1422.49482387392"
$ws.Cells.Item(108,5).Value = "'1460.89789759927"
$ws.Cells.Item(109,5).Value = "'1593.89826707408"
$ws.Cells.Item(110,5).Value = "'1561.02176377286"
$ws.Cells.Item(111,5).Value = "'1596.18813913168"
$ws.Cells.Item(112,5).Value = "'1725.59684609893"
$ws.Cells.Item(113,5).Value = "'1831.31393307957"
$ws.Cells.Item(114,5).Value = "'1837.19395628138"
$ws.Cells.Item(115,5).Value = "'1895.45328192576"
$ws.Cells.Item(116,5).Value = "'1793.02814065952"
$ws.Cells.Item(117,5).Value = "'1996.16491123244"
$ws.Cells.Item(118,5).Value = "'2121.92475626296"
$ws.Cells.Item(119,5).Value = "'2469.67781941679"
$ws.Cells.Item(120,5).Value = "'2557.07319722189"
$ws.Cells.Item(121,5).Value = "'2494.33784671389"
$ws.Cells.Item(122,5).Value = "'2623.73050377988"
$ws.Cells.Item(123,5).Value = "'2621.04329660168"
$ws.Cells.Item(124,5).Value = "'2399.6421784848"
$ws.Cells.Item(125,5).Value = "'2462.00463656339"
$ws.Cells.Item(132,5).Value = "'"
$ws.Cells.Item(133,5).Value = "'"
$ws.Cells.Item(134,5).Value = "'"
$ws.Cells.Item(135,5).Value = "'"
$ws.Cells.Item(136,5).Value = "'"
$ws.Cells.Item(137,5).Value = "'"
$ws.Cells.Item(138,5).Value = "'"
$ws.Cells.Item(139,5).Value = "'"
$ws.Cells.Item(140,5).Value = "'"
$ws.Cells.Item(141,5).Value = "'"
$ws.Cells.Item(142,5).Value = "'"
$ws.Cells.Item(143,5).Value = "'"
$ws.Cells.Item(144,5).Value = "'"
$ws.Cells.Item(145,5).Value = "'"
$ws.Cells.Item(146,5).Value = "'"
$ws.Cells.Item(147,5).Value = "'"
$ws.Cells.Item(148,5).Value = "'"
$ws.Cells.Item(149,5).Value = "'"
$ws.Cells.Item(150,5).Value = "'"
$ws.Cells.Item(151,5).Value = "'"
$ws.Cells.Item(152,5).Value = "'"
$ws.Cells.Item(153,5).Value = "'"
$ws.Cells.Item(154,5).Value = "'"
$ws.Cells.Item(155,5).Value = "'"
$ws.Cells.Item(156,5).Value = "'"
$ws.Cells.Item(157,5).Value = "'"
$ws.Cells.Item(158,5).Value = "'"
$ws.Cells.Item(159,5).Value = "'"
$ws.Cells.Item(160,5).Value = "'"
$ws.Cells.Item(161,5).Value = "'"
$ws.Cells.Item(162,5).Value = "'"
$ws.Cells.Item(163,5).Value = "'"
$ws.Cells.Item(164,5).Value = "'"
$ws.Cells.Item(165,5).Value = "'"
$ws.Cells.Item(166,5).Value = "'"
$ws.Cells.Item(167,5).Value = "'"
$ws.Cells.Item(168,5).Value = "'"
$ws.Cells.Item(169,5).Value = "'"
$ws.Cells.Item(170,5).Value = "'"
$ws.Cells.Item(171,5).Value = "'"
$ws.Cells.Item(172,5).Value = "'2455.81857694823"
$ws.Cells.Item(173,5).Value = "'2315.95886869544"
$ws.Cells.Item(174,5).Value = "'2120.54612462779"
$ws.Cells.Item(175,5).Value = "'1996.29334047591"
$ws.Cells.Item(176,5).Value = "'1925.9560329888"
$ws.Cells.Item(177,5).Value = "'1816.00900994004"
$ws.Cells.Item(178,5).Value = "'1733.09950574924"
$ws.Cells.Item(179,5).Value = "'1604.89952674453"
$ws.Cells.Item(180,5).Value = "'1579.93092964188"
$ws.Cells.Item(181,5).Value = "'1664.64215028779"
$ws.Cells.Item(182,5).Value = "'1657.14443577225"
$ws.Cells.Item(183,5).Value = "'1704.37326161717"
$ws.Cells.Item(184,5).Value = "'1710.94957201576"
$ws.Cells.Item(185,5).Value = "'1729.62207872068"
$ws.Cells.Item(186,5).Value = "'1753.69193369254"
$ws.Cells.Item(187,5).Value = "'1808.06182579764"
$ws.Cells.Item(188,5).Value = "'1778.49565915524"
$ws.Cells.Item(189,5).Value = "'1746.92756249456"
$ws.Cells.Item(190,5).Value = "'1794.65737167766"
$ws.Cells.Item(191,1).Value = 408
$ws.Cells.Item(191,2).Value = "North Korea"
$ws.Cells.Item(191,3).Value = "GDP per Capita"
$ws.Cells.Item(191,4).Value = 2009
$ws.Cells.Item(191,5).Value = "'1768.88974979975"
$ws.Cells.Item(192,1).Value = 408
$ws.Cells.Item(192,2).Value = "North Korea"
$ws.Cells.Item(192,3).Value = "GDP per Capita"
$ws.Cells.Item(192,4).Value = 2010
$ws.Cells.Item(192,5).Value = "'1751.43628161136"
$ws.Cells.Item(193,1).Value = 408
$ws.Cells.Item(193,2).Value = "North Korea"
$ws.Cells.Item(193,3).Value = "GDP per Capita"
$ws.Cells.Item(193,4).Value = 2011
$ws.Cells.Item(193,5).Value = "'1756.61080668874"
$ws.Cells.Item(194,1).Value = 408
$ws.Cells.Item(194,2).Value = "North Korea"
$ws.Cells.Item(194,3).Value = "GDP per Capita"
$ws.Cells.Item(194,4).Value = 2012
$ws.Cells.Item(194,5).Value = "'1770.92094154285"
$ws.Cells.Item(195,1).Value = 408
$ws.Cells.Item(195,2).Value = "North Korea"
$ws.Cells.Item(195,3).Value = "GDP per Capita"
$ws.Cells.Item(195,4).Value = 2013
$ws.Cells.Item(195,5).Value = "'1781.30084534494"
$ws.Cells.Item(196,1).Value = 408
$ws.Cells.Item(196,2).Value = "North Korea"
$ws.Cells.Item(196,3).Value = "GDP per Capita"
$ws.Cells.Item(196,4).Value = 2014
$ws.Cells.Item(196,5).Value = "'1791.34364009763"
$ws.Cells.Item(197,1).Value = 408
$ws.Cells.Item(197,2).Value = "North Korea"
$ws.Cells.Item(197,3).Value = "GDP per Capita"
$ws.Cells.Item(197,4).Value = 2015
$ws.Cells.Item(197,5).Value = "'1762.51667533842"
$ws.Cells.Item(198,1).Value = 408
$ws.Cells.Item(198,2).Value = "North Korea"
$ws.Cells.Item(198,3).Value = "GDP per Capita"
$ws.Cells.Item(198,4).Value = 2016
$ws.Cells.Item(198,5).Value = "'1742.02229539263"
